$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row of coverage history data (row 6) following the same
# pattern as the existing rows (date in column A, raw counts in C:K,
# computed percentage formulas in M:Q).

# Copy the date format from A5 to A6 (keeps the same underlying style,
# rather than manufacturing a brand new custom number format).
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Copy the 2-decimal percentage number format from M5:Q5 to M6:Q6.
$ws.Range("M5:Q5").Copy()
$ws.Range("M6:Q6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A6").Value = 44058

$ws.Range("C6").Value = 161
$ws.Range("D6").Value = 73
$ws.Range("E6").Value = 1445
$ws.Range("F6").Value = 439
$ws.Range("G6").Value = 35
$ws.Range("H6").Value = 20
$ws.Range("I6").Value = 2
$ws.Range("J6").Value = 209
$ws.Range("K6").Value = 83

$ws.Range("M6").Formula = "=100*D6/C6"
$ws.Range("N6").Formula = "=100*F6/E6"
$ws.Range("O6").Formula = "=100*G6/C6"
$ws.Range("P6").Formula = "=100*I6/H6"
$ws.Range("Q6").Formula = "=100*K6/J6"

# Move selection to the next empty row, as Excel would after data entry.
$ws.Range("A7").Select()

$wb.Save()
